$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.316.99"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.877.42"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "'242.43"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.08025"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("D9").Value = "'0.3170"
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("D10").Value = "'24.98"
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("D11").Value = "'0.08312"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").Value = "1.887.48"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "'94.64"
$ws.Range("E14").Value = "  +3.73%  "
$ws.Range("D15").Value = "'0.7162"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "'6.394"
$ws.Range("E16").Value = "  +4.84%  "
$ws.Range("D17").Value = "'0.000008639"
$ws.Range("E17").Value = "  +4.80%  "
$ws.Range("D18").Value = "29.351.27"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "'242.73"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").Value = "'13.32"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").Value = "2.139.77"
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "'0.1568"
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("D26").Value = "'9.086"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "'162.93"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "'18.55"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'1.506"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "'4.424"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'4.333"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "'1.191"
$ws.Range("E32").Value = "  -7.78%  "
$ws.Range("D33").Value = "'0.05404"
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("D34").Value = "'1.942"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "'0.7728"
$ws.Range("E35").Value = "  +3.89%  "
$ws.Range("D36").Value = "'1.187"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").Value = "'2.683"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").Value = "'0.01887"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").Value = "1.264.55"
$ws.Range("E39").Value = "  +3.07%  "
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").Value = "'6.492"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").Value = "'113.26"
$ws.Range("E42").Value = "  +2.22%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.9071"
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'74.49"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("E45").Value = "  +7.83%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "2.029.87"
$ws.Range("D48").Value = "'1.809"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").Value = "'0.5222"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").Value = "'9.503"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").Value = "'0.4371"
$ws.Range("E51").Value = "  +1.12%  "
